# Social Cost of Carbon.xlsx - content update
#
# Summary of the edit (from commit "Updated Indonesia files compatible with v3.3.1"):
#  - About sheet: the "Social Cost of Carbon" note sentence is split across two
#    lines and reworded; the CPI conversion-factor value (1.109) shifts down one
#    row (a new row is inserted above it).
#  - SourceData sheet: every H:K formula that referenced the CPI factor at
#    About!$A$16 now points at About!$A$17 (because of the inserted row above).
#  - SCoC sheet: the header label for the $ column is reworded to include units.
#  - Selections on each sheet are left the way the author's Excel session had
#    them when saved.

$wb = $excel.ActiveWorkbook

$about      = $wb.Worksheets.Item("About")
$sourceData = $wb.Worksheets.Item("SourceData")
$scoc       = $wb.Worksheets.Item("SCoC")

# ---------------------------------------------------------------------------
# About sheet: insert a new row above the old "a 3% discount rate..." row
# (row 11) so the long intro sentence can be split in two, and reword it.
# ---------------------------------------------------------------------------
$about.Rows.Item(11).Insert()

$about.Range("A10").Value = "When considering the Social Cost of Carbon, meant to capture the long-term economic damage caused by one"
$about.Range("A11").Value = "ton of carbon dioxide emitted, the U.S. government typically uses the figures based on"

# ---------------------------------------------------------------------------
# SourceData sheet: the CPI conversion factor cell used to live at About!A16;
# after the row insert above it now lives at About!A17, so repoint every
# formula in the H:K calculation block (rows 4-44).
# ---------------------------------------------------------------------------
$null = $sourceData.Range("H4:K44").Replace("About!`$A`$16", "About!`$A`$17")

# ---------------------------------------------------------------------------
# SCoC sheet: clarify the units in the column header.
# ---------------------------------------------------------------------------
$scoc.Range("B1").Value = "Social Cost of Carbon (`$/g CO2e)"

# ---------------------------------------------------------------------------
# Restore the per-sheet selections / active cell the author's session had,
# finishing on the About sheet so it stays the active tab.
# ---------------------------------------------------------------------------
$null = $sourceData.Activate()
$null = $sourceData.Range("B2").Select()

$null = $scoc.Activate()
$null = $scoc.Range("B2").Select()

$null = $about.Activate()
$null = $about.Range("D10").Select()
